$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Record "Completed Effort" for the two newly-finished features on Day 3 and
# Day 4 (columns G and H), rows 13 and 14.
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1

# Fix the "Ideal Burndown" formula/style for Day 1 (E20) so it follows the
# same shared pattern as the rest of the row instead of the stray
# "D20-E19" formula. Copy F20's formatting (fill/border/number format) onto
# E20 first, then give it the proper ideal-burndown formula.
$ws.Range("F20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Formula = '=$D$20-($D$20/$K$5*E5)'

# Match the author's final selection when the file was saved.
$ws.Range("H14").Select()
